# Commit: "changing names and activating google sheet"
#  - delete Sheet4 (duplicate data sheet, no longer needed)
#  - normalise regStatus labels to lower-case ("removed" / "unknown")
#  - move the active selection on Sheet1 to K12

$wb = $excel.ActiveWorkbook

# --- remove the now-unused Sheet4 -----------------------------------------
$excel.DisplayAlerts = $false
$wb.Worksheets.Item("Sheet4").Delete() | Out-Null
$excel.DisplayAlerts = $true

$ws1 = $wb.Worksheets.Item("Sheet1")

# --- normalise the regStatus column (K) on Sheet1 --------------------------
# Rows 3-11 all become "removed"; row 12 becomes "unknown".
$ws1.Range("K3").Value = "removed"
$ws1.Range("K4").Value = "removed"
$ws1.Range("K5").Value = "removed"
$ws1.Range("K6").Value = "removed"
$ws1.Range("K7").Value = "removed"
$ws1.Range("K8").Value = "removed"
$ws1.Range("K9").Value = "removed"
$ws1.Range("K10").Value = "removed"
$ws1.Range("K11").Value = "removed"
$ws1.Range("K12").Value = "unknown"

# --- update the view: activate Sheet1 and move the selection to K12 -------
$ws1.Activate() | Out-Null
$ws1.Range("K12").Select() | Out-Null
